$wb = $excel.ActiveWorkbook

# Sheet pojos_1
$ws1 = $wb.Worksheets.Item("pojos_1")
$ws1.Range("C2").Value = 889182.0
$ws1.Range("D2").Value = 914770.75
$ws1.Range("E2").Value = 188138.0
$ws1.Range("C3").Value = 203702.0
$ws1.Range("D3").Value = 754265.3125
$ws1.Range("E3").Value = 861124.0

# Sheet pojos_2
$ws2 = $wb.Worksheets.Item("pojos_2")
$ws2.Range("C2").Value = 270245.0
$ws2.Range("D2").Value = 533421.75
$ws2.Range("E2").Value = 123410.0
$ws2.Range("C3").Value = 851954.0
$ws2.Range("D3").Value = 344267.78125
$ws2.Range("E3").Value = 298359.0

# Sheet pojos_3
$ws3 = $wb.Worksheets.Item("pojos_3")
$ws3.Range("C2").Value = 222746.0
$ws3.Range("D2").Value = 116392.609375
$ws3.Range("E2").Value = 347667.0
